$d = $word.ActiveDocument

# 1) Merge the split runs "reproduzir: " + text into a single run.
$r1 = $d.Content
$r1.Find.Execute(": cria 2 novas", $true, $false, $false, $false, $false, `
                  $true, 1, $false, ": cria 2 novas", 2) | Out-Null

# 2) Merge the split runs "obterValor: " + text into a single run.
$r2 = $d.Content
$r2.Find.Execute(": obtém o valor", $true, $false, $false, $false, $false, `
                  $true, 1, $false, ": obtém o valor", 2) | Out-Null

# 3) Move the automatic "_GoBack" bookmark from the end of the document
#    (after the last author's e-mail address) to the empty centered
#    paragraph that immediately precedes the title - this is where Word
#    last left the edit point.
$target = $d.Paragraphs.Item(9)
$d.Bookmarks.Add("_GoBack", $target.Range) | Out-Null
